$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'97.678.29"
$ws.Range("E2").Value = "'  -1.76%  "
$ws.Range("D3").Value = "'3.432.14"
$ws.Range("E3").Value = "'  +4.13%  "
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'256.29"
$ws.Range("E5").Value = "'  +0.94%  "
$ws.Range("D6").Value = "'658.00"
$ws.Range("E6").Value = "'  +5.52%  "
$ws.Range("E7").Value = "'  +5.45%  "
$ws.Range("D8").Value = "'0.431"
$ws.Range("E8").Value = "'  +7.86%  "
$ws.Range("E9").Value = "'  +10.38%  "
$ws.Range("E10").Value = "'  +0.02%  "
$ws.Range("D11").Value = "'3.429.39"
$ws.Range("E11").Value = "'  +4.09%  "
$ws.Range("E12").Value = "'  +6.96%  "
$ws.Range("D13").Value = "'42.28"
$ws.Range("E13").Value = "'  +6.72%  "
$ws.Range("D14").Value = "'6.39"
$ws.Range("E14").Value = "'  +16.63%  "
$ws.Range("E15").Value = "'  +5.74%  "
$ws.Range("D16").Value = "'97.359.42"
$ws.Range("E16").Value = "'  -1.78%  "
$ws.Range("D17").Value = "'4.063.46"
$ws.Range("E17").Value = "'  +4.77%  "
$ws.Range("D18").Value = "'8.76"
$ws.Range("E18").Value = "'  +37.66%  "
$ws.Range("D19").Value = "'3.424.15"
$ws.Range("E19").Value = "'  +4.34%  "
$ws.Range("D20").Value = "'17.73"
$ws.Range("E20").Value = "'  +14.79%  "
$ws.Range("D21").Value = "'0.523"
$ws.Range("E21").Value = "'  +62.46%  "
$ws.Range("D22").Value = "'11.00"
$ws.Range("E22").Value = "'  +17.74%  "
$ws.Range("E23").Value = "'  +1.34%  "
$ws.Range("D24").Value = "'510.92"
$ws.Range("E24").Value = "'  +4.21%  "
$ws.Range("D25").Value = "'0.0000208"
$ws.Range("E25").Value = "'  +3.16%  "
$ws.Range("D26").Value = "'6.23"
$ws.Range("E26").Value = "'  +10.14%  "
$ws.Range("D27").Value = "'99.85"
$ws.Range("E27").Value = "'  +12.28%  "
$ws.Range("D28").Value = "'12.90"
$ws.Range("E28").Value = "'  +7.22%  "
$ws.Range("D29").Value = "'3.609.20"
$ws.Range("E29").Value = "'  +5.25%  "
$ws.Range("D30").Value = "'0.157"
$ws.Range("E30").Value = "'  +14.67%  "
$ws.Range("D31").Value = "'11.52"
$ws.Range("E31").Value = "'  +10.93%  "
$ws.Range("D32").Value = "'0.200"
$ws.Range("E32").Value = "'  +5.68%  "
$ws.Range("E33").Value = "'  -0.09%  "
$ws.Range("E34").Value = "'  +22.45%  "
$ws.Range("E35").Value = "'  +0.39%  "
$ws.Range("E36").Value = "'  +7.75%  "
$ws.Range("D37").Value = "'2.30"
$ws.Range("E37").Value = "'  +18.38%  "
$ws.Range("D38").Value = "'7.89"
$ws.Range("E38").Value = "'  +9.08%  "
$ws.Range("E39").Value = "'  +4.44%  "
$ws.Range("D40").Value = "'1.44"
$ws.Range("E40").Value = "'  +16.77%  "
$ws.Range("D41").Value = "'521.37"
$ws.Range("E41").Value = "'  +6.47%  "
$ws.Range("D42").Value = "'24.72"
$ws.Range("E42").Value = "'  -0.42%  "
$ws.Range("D43").Value = "'0.871"
$ws.Range("E43").Value = "'  +12.52%  "
$ws.Range("D44").Value = "'3.69"
$ws.Range("E44").Value = "'  +1.70%  "
$ws.Range("D45").Value = "'0.0423"
$ws.Range("E45").Value = "'  +26.12%  "
$ws.Range("D46").Value = "'3.33"
$ws.Range("E46").Value = "'  +7.49%  "
$ws.Range("D47").Value = "'5.54"
$ws.Range("E47").Value = "'  +16.80%  "
$ws.Range("D48").Value = "'8.28"
$ws.Range("E48").Value = "'  +13.30%  "
$ws.Range("E49").Value = "'  +0.08%  "
$ws.Range("E50").Value = "'  +17.28%  "
$ws.Range("E51").Value = "'  +7.31%  "
